$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header "ecc" in C1 (adds new shared string)
$ws.Range("C1").Value = "ecc"

# C2: first occurrence of the formula, entered standalone (not shared)
$ws.Range("C2").Formula = "=IF(ABS(A2)>=ABS(B2),ABS(A2),ABS(B2))"

# C3:C66 filled together as one shared-formula group
$ws.Range("C3:C66").Formula = "=IF(ABS(A3)>=ABS(B3),ABS(A3),ABS(B3))"

# C67:C69 filled together as a second shared-formula group
$ws.Range("C67:C69").Formula = "=IF(ABS(A67)>=ABS(B67),ABS(A67),ABS(B67))"

# Update the active selection as recorded in the saved workbook
$ws.Range("F70").Select() | Out-Null
